$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row for the "Prev 30-day Precip" layer, right above the
# --- "Fire Intensity Risk" row (old row 19), shifting everything below down.
$ws.Rows.Item(19).Insert()

# --- Populate the new row 19 with the 30-day precipitation layer info ---
$ws.Range("A19").Value = "Prev 30-day Precip"
$ws.Range("B19").Value = "30-day total precipitation"
$ws.Range("C19").Value = "USGS"
$ws.Range("D19").Value = "WMS Image"
$ws.Range("F19").Value = "PRECIP_TP30"

# E19 carries a hyperlink (same target as the "Prev 7-day Precip" row above it)
$ws.Range("E19").Value = "https://vegdri.cr.usgs.gov/wms.php?"
$ws.Hyperlinks.Add($ws.Range("E19"), "https://vegdri.cr.usgs.gov/wms.php?") | Out-Null
# Re-apply the plain hyperlink style (matches the style used by E18) since
# Hyperlinks.Add() on this range can otherwise mint a duplicate style entry.
$ws.Range("E19").Style = $ws.Range("E18").Style

# --- Restore the view state: scrolled down a bit with F19 as the active cell ---
$ws.Range("F19").Select()
